{"js": "// Replace the date line.\nconst dateResults = context.document.body.search(\"2024-06-18 Tuesday\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\ndateResults.items[0].insertText(\"2024-06-19 Wednesday\", Word.InsertLocation.replace);\n\n// Replace each division expression in the practice table. Values that are\n// unique in the document are searched directly; \"24\u00f79=\" occurs twice, so it\n// is searched once and its two matches are replaced in document order.\nconst singleReplacements = [\n  [\"68\u00f73=\", \"72\u00f75=\"],\n  [\"96\u00f78=\", \"53\u00f73=\"],\n  [\"65\u00f76=\", \"46\u00f78=\"],\n  [\"89\u00f77=\", \"74\u00f73=\"],\n  [\"31\u00f79=\", \"77\u00f77=\"],\n  [\"19\u00f73=\", \"82\u00f74=\"],\n  [\"21\u00f73=\", \"33\u00f78=\"],\n  [\"11\u00f72=\", \"71\u00f76=\"],\n  [\"45\u00f76=\", \"28\u00f75=\"],\n  [\"95\u00f72=\", \"66\u00f74=\"],\n  [\"78\u00f79=\", \"10\u00f78=\"],\n  [\"79\u00f74=\", \"81\u00f72=\"],\n  [\"19\u00f79=\", \"28\u00f77=\"],\n  [\"76\u00f72=\", \"36\u00f72=\"],\n  [\"50\u00f79=\", \"92\u00f77=\"],\n  [\"62\u00f75=\", \"80\u00f73=\"],\n  [\"97\u00f78=\", \"38\u00f77=\"],\n  [\"43\u00f73=\", \"20\u00f72=\"],\n  [\"69\u00f73=\", \"42\u00f73=\"],\n  [\"57\u00f75=\", \"91\u00f73=\"],\n  [\"13\u00f74=\", \"12\u00f76=\"],\n  [\"86\u00f76=\", \"92\u00f75=\"],\n  [\"55\u00f72=\", \"12\u00f77=\"],\n];\n\nfor (const [oldText, newText] of singleReplacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\n// \"24\u00f79=\" appears twice: first occurrence -> \"41\u00f77=\", second -> \"15\u00f72=\".\nconst dupResults = context.document.body.search(\"24\u00f79=\", { matchCase: true });\ndupResults.load(\"items\");\nawait context.sync();\ndupResults.items[0].insertText(\"41\u00f77=\", Word.InsertLocation.replace);\ndupResults.items[1].insertText(\"15\u00f72=\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Update the heading date, then update each division expression in the\n# practice table. Word constants used directly (no $Word enum available):\n#   wdReplaceOne = 1, wdReplaceAll = 2, wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\nfunction Replace-OneText($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1) | Out-Null\n}\n\n# Heading date.\nReplace-AllText \"2024-06-18 Tuesday\" \"2024-06-19 Wednesday\"\n\n# Division expressions that are unique in the document.\nReplace-AllText \"68\u00f73=\" \"72\u00f75=\"\nReplace-AllText \"96\u00f78=\" \"53\u00f73=\"\nReplace-AllText \"65\u00f76=\" \"46\u00f78=\"\nReplace-AllText \"89\u00f77=\" \"74\u00f73=\"\nReplace-AllText \"31\u00f79=\" \"77\u00f77=\"\nReplace-AllText \"19\u00f73=\" \"82\u00f74=\"\nReplace-AllText \"21\u00f73=\" \"33\u00f78=\"\nReplace-AllText \"11\u00f72=\" \"71\u00f76=\"\nReplace-AllText \"45\u00f76=\" \"28\u00f75=\"\nReplace-AllText \"95\u00f72=\" \"66\u00f74=\"\nReplace-AllText \"78\u00f79=\" \"10\u00f78=\"\nReplace-AllText \"79\u00f74=\" \"81\u00f72=\"\nReplace-AllText \"19\u00f79=\" \"28\u00f77=\"\nReplace-AllText \"76\u00f72=\" \"36\u00f72=\"\nReplace-AllText \"50\u00f79=\" \"92\u00f77=\"\nReplace-AllText \"62\u00f75=\" \"80\u00f73=\"\nReplace-AllText \"97\u00f78=\" \"38\u00f77=\"\nReplace-AllText \"43\u00f73=\" \"20\u00f72=\"\nReplace-AllText \"69\u00f73=\" \"42\u00f73=\"\nReplace-AllText \"57\u00f75=\" \"91\u00f73=\"\nReplace-AllText \"13\u00f74=\" \"12\u00f76=\"\nReplace-AllText \"86\u00f76=\" \"92\u00f75=\"\nReplace-AllText \"55\u00f72=\" \"12\u00f77=\"\n\n# \"24\u00f79=\" appears twice in document order: the first cell becomes \"41\u00f77=\",\n# the second becomes \"15\u00f72=\". Replace one occurrence at a time so each\n# instance gets its own target text.\nReplace-OneText \"24\u00f79=\" \"41\u00f77=\"\nReplace-OneText \"24\u00f79=\" \"15\u00f72=\"\n"}
